$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Tim David"

# Insert a new column before column A, shifting existing data to the right
$ws.Columns("A").Insert()

# Fill in the new column A values
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "35th"
